$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for each data row (2-216).
# Update every occurrence of the old date serial (45172) to the new one (45175),
# preserving the existing date number format / style on the cells.
$ws.Range("C2:C216").Value2 = 45175
